$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates
$ws.Range("B1").Value = "is.global"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "latitude"
$ws.Range("E1").Value = "longitude"

# Data row updates
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Update selection to match target (E3)
$ws.Range("E3").Select()
